$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("apiTest")

$ws2.Range("E1").Value = "newFirstName"
$ws2.Range("F1").Value = "newLastName"
$ws2.Range("G1").Value = "newPass"

$ws2.Range("E2").Value = "Srdjan1"
$ws2.Range("F2").Value = "Rados1"
$ws2.Range("G2").Value = "Test123@"

$ws2.Activate() | Out-Null
$ws2.Range("G2").Select() | Out-Null
